$d = $word.ActiveDocument

# 1. Replace the placeholder hosting-URL text with the actual video link.
$d.Content.Find.Execute("The Hosting URL that I can use to access your application.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://www.youtube.com/watch?v=THb9KMIPtjE", 2)

# 2. Remove the old "_GoBack" bookmark (it used to sit after the "Y" in the
#    Peer Review row - Word always keeps only one "_GoBack" bookmark, moving
#    it to the location of the most recent edit).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Re-create "_GoBack" right after the newly inserted URL text (collapsed,
#    zero-length bookmark), matching Word's behaviour of tracking the last
#    edit location. A truly collapsed range positioned exactly at a run/
#    paragraph-mark boundary cannot be fed straight into Bookmarks.Add, so
#    build it via a tiny placeholder run that is removed immediately after.
$urlRange = $d.Content
$urlRange.Find.Execute("https://www.youtube.com/watch?v=THb9KMIPtjE")
$urlRange.Collapse(0)
$urlRange.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $urlRange)
$goBack = $d.Bookmarks("_GoBack")
$goBack.Range.Text = ""
